# Completed Data visualizing using tableau and added dashboards
# Adds the "Asia" continent row (row 8) to the COVID world-analysis table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row: Asia ---
# A8: Location
$ws.Range("A8").Value = "Asia"

# B8: population (typed with the sheet's normal font explicitly applied)
$ws.Range("B8").Value = 4721383370
$ws.Range("B8").Font.Name = "Aptos Narrow"

# C8: total_case
$ws.Range("C8").Value = 301428726

# D8: total_death
$ws.Range("D8").Value = 1637060

# E8: infection_rate = (total_case / population) * 100
$ws.Range("E8").Formula = "=(C8/B8)*100"

# F8: death_rate = (total_death / population) * 100
$ws.Range("F8").Formula = "=(D8/B8)*100"

# Both calculated columns are formatted with 2 decimal places
$ws.Range("E8:F8").NumberFormat = "0.00"

# Move the active selection the way it ended up after entering the row
$ws.Range("F9").Select()
